# Update "想去人数" (F column) counts across the workbook's four sheets.
# Sheet1 "展览" (Exhibitions), Sheet2 "演出" (Performances), Sheet3 "本地生活"
# (Local Life) hold the per-category rows; Sheet4 "全部类型" (All Types) is an
# aggregated view of the same events, so matching rows are updated there too.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1)
$ws1.Cells.Item(3, 6).Value = 1405
$ws1.Cells.Item(4, 6).Value = 25630
$ws1.Cells.Item(6, 6).Value = 238
$ws1.Cells.Item(8, 6).Value = 160
$ws1.Cells.Item(12, 6).Value = 195
$ws1.Cells.Item(13, 6).Value = 168
$ws1.Cells.Item(15, 6).Value = 268
$ws1.Cells.Item(16, 6).Value = 328
$ws1.Cells.Item(17, 6).Value = 44
$ws1.Cells.Item(18, 6).Value = 1462
$ws1.Cells.Item(19, 6).Value = 149
$ws1.Cells.Item(20, 6).Value = 405
$ws1.Cells.Item(21, 6).Value = 89

# Sheet "演出" (sheet2)
$ws2.Cells.Item(6, 6).Value = 64
$ws2.Cells.Item(14, 6).Value = 14

# Sheet "本地生活" (sheet3)
$ws3.Cells.Item(2, 6).Value = 4865

# Sheet "全部类型" (sheet4) - aggregated rows for the same events
$ws4.Cells.Item(3, 6).Value = 1405
$ws4.Cells.Item(4, 6).Value = 4865
$ws4.Cells.Item(6, 6).Value = 25630
$ws4.Cells.Item(9, 6).Value = 238
$ws4.Cells.Item(14, 6).Value = 160
$ws4.Cells.Item(15, 6).Value = 64
$ws4.Cells.Item(16, 6).Value = 64
$ws4.Cells.Item(25, 6).Value = 195
$ws4.Cells.Item(26, 6).Value = 168
$ws4.Cells.Item(29, 6).Value = 268
$ws4.Cells.Item(31, 6).Value = 14
$ws4.Cells.Item(32, 6).Value = 328
$ws4.Cells.Item(33, 6).Value = 44
$ws4.Cells.Item(35, 6).Value = 1462
$ws4.Cells.Item(36, 6).Value = 149
$ws4.Cells.Item(38, 6).Value = 405
$ws4.Cells.Item(39, 6).Value = 89
